$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.243623333333334
$ws.Range("H2").Value = 12.73087
$ws.Range("I2").Value = 0.2469246453968972
$ws.Range("J2").Value = 0.2469246453968973
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.328871
$ws.Range("N2").Value = 6.986613
$ws.Range("O2").Value = 0.2411747569970185
$ws.Range("P2").Value = 0.2411747569970185
$ws.Range("Q2").Value = 9.882851315923334
$ws.Range("R2").Value = 88.94566184331
$ws.Range("S2").Value = 0.05955199135017164
$ws.Range("T2").Value = 0.05955199135017166

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.243623333333334
$ws.Range("H3").Value = 12.73087
$ws.Range("I3").Value = 0.2469246453968972
$ws.Range("J3").Value = 0.2469246453968973
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.327491999999999
$ws.Range("N3").Value = 21.982476
$ws.Range("O3").Value = 0.7588252430029816
$ws.Range("P3").Value = 0.7588252430029816
$ws.Range("Q3").Value = 31.09511602601333
$ws.Range("R3").Value = 279.85604423412
$ws.Range("S3").Value = 0.1873726540467256
$ws.Range("T3").Value = 0.1873726540467257

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.865491666666667
$ws.Range("H4").Value = 17.596475
$ws.Range("I4").Value = 0.3412966552647515
$ws.Range("J4").Value = 0.3412966552647516
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.328871
$ws.Range("N4").Value = 6.986613
$ws.Range("O4").Value = 0.2411747569970185
$ws.Range("P4").Value = 0.2411747569970185
$ws.Range("Q4").Value = 13.65997344324167
$ws.Range("R4").Value = 122.939760989175
$ws.Range("S4").Value = 0.08231213789737163
$ws.Range("T4").Value = 0.08231213789737166

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.865491666666667
$ws.Range("H5").Value = 17.596475
$ws.Range("I5").Value = 0.3412966552647515
$ws.Range("J5").Value = 0.3412966552647516
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.327491999999999
$ws.Range("N5").Value = 21.982476
$ws.Range("O5").Value = 0.7588252430029816
$ws.Range("P5").Value = 0.7588252430029816
$ws.Range("Q5").Value = 42.97934326356667
$ws.Range("R5").Value = 386.8140893721
$ws.Range("S5").Value = 0.2589845173673799
$ws.Range("T5").Value = 0.25898451736738

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.123111999999999
$ws.Range("H6").Value = 18.369336
$ws.Range("I6").Value = 0.3562868663317164
$ws.Range("J6").Value = 0.3562868663317164
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.328871
$ws.Range("N6").Value = 6.986613
$ws.Range("O6").Value = 0.2411747569970185
$ws.Range("P6").Value = 0.2411747569970185
$ws.Range("Q6").Value = 14.259937966552
$ws.Range("R6").Value = 128.339441698968
$ws.Range("S6").Value = 0.08592739840878089
$ws.Range("T6").Value = 0.08592739840878091

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.123111999999999
$ws.Range("H7").Value = 18.369336
$ws.Range("I7").Value = 0.3562868663317164
$ws.Range("J7").Value = 0.3562868663317164
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.327491999999999
$ws.Range("N7").Value = 21.982476
$ws.Range("O7").Value = 0.7588252430029816
$ws.Range("P7").Value = 0.7588252430029816
$ws.Range("Q7").Value = 44.86705419510399
$ws.Range("R7").Value = 403.8034877559359
$ws.Range("S7").Value = 0.2703594679229355
$ws.Range("T7").Value = 0.2703594679229355

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.9536773333333334
$ws.Range("H8").Value = 2.861032
$ws.Range("I8").Value = 0.05549183300663471
$ws.Range("J8").Value = 0.05549183300663472
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.328871
$ws.Range("N8").Value = 6.986613
$ws.Range("O8").Value = 0.2411747569970185
$ws.Range("P8").Value = 0.2411747569970185
$ws.Range("Q8").Value = 2.220991484957334
$ws.Range("R8").Value = 19.988923364616
$ws.Range("S8").Value = 0.01338322934069425
$ws.Range("T8").Value = 0.01338322934069426

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.9536773333333334
$ws.Range("H9").Value = 2.861032
$ws.Range("I9").Value = 0.05549183300663471
$ws.Range("J9").Value = 0.05549183300663472
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.327491999999999
$ws.Range("N9").Value = 21.982476
$ws.Range("O9").Value = 0.7588252430029816
$ws.Range("P9").Value = 0.7588252430029816
$ws.Range("Q9").Value = 6.988063030581333
$ws.Range("R9").Value = 62.892567275232
$ws.Range("S9").Value = 0.04210860366594046
$ws.Range("T9").Value = 0.04210860366594047
